# Applies the "Updated cryptos list" data refresh (commit: Sat Jun 15 23:59:04 UTC 2024,
# GitHub Actions) to the cryptos worksheet.
#
# - Price (column D) / Volume(1h) (column E) numbers are refreshed for most rows.
# - Three coin pairs were re-ranked and swapped rows (33/34, 40/41, 50/51), so their
#   Coin name (B), Link (C), Price (D) and Volume (E) are rewritten in place.
#
# All D/E values in the source sheet are stored as *text* (t="inlineStr"), including
# price strings that look numeric (e.g. "0.997", "608.34"). Excel's COM layer will
# silently coerce a plain numeric-looking string assigned to .Value into a real
# number, so we force the Text number format on the data range first to keep
# every value a string, matching the original file's cell typing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '66.234.77'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '3.566.63'
$ws.Range("E3").Value = '  +2.43%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '608.34'
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("D6").Value = '145.46'
$ws.Range("E6").Value = '  +1.64%  '
$ws.Range("D7").Value = '3.565.77'
$ws.Range("E7").Value = '  +2.43%  '
$ws.Range("E9").Value = '  +3.48%  '
$ws.Range("E10").Value = '  +1.16%  '
$ws.Range("D11").Value = '7.91'
$ws.Range("E11").Value = '  -3.24%  '
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").Value = '4.171.35'
$ws.Range("E13").Value = '  +2.49%  '
$ws.Range("E14").Value = '  +2.37%  '
$ws.Range("D15").Value = '30.00'
$ws.Range("E15").Value = '  -1.01%  '
$ws.Range("D16").Value = '3.577.02'
$ws.Range("E16").Value = '  +2.72%  '
$ws.Range("D17").Value = '66.329.04'
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("D19").Value = '11.51'
$ws.Range("E19").Value = '  +10.38%  '
$ws.Range("D20").Value = '6.22'
$ws.Range("E20").Value = '  +0.93%  '
$ws.Range("E21").Value = '  +0.74%  '
$ws.Range("D22").Value = '429.92'
$ws.Range("E22").Value = '  +1.61%  '
$ws.Range("E23").Value = '  +4.43%  '
$ws.Range("D25").Value = '3.709.27'
$ws.Range("E25").Value = '  +2.61%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  +3.64%  '
$ws.Range("E28").Value = '  +2.21%  '
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("E30").Value = '  -2.41%  '
$ws.Range("D31").Value = '0.997'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("D32").Value = '25.63'
$ws.Range("E32").Value = '  +1.94%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '1.46'
$ws.Range("E33").Value = '  -1.48%  '
$ws.Range("B34").Value = 'RenzoRestakedETH'
$ws.Range("C34").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D34").Value = '3.563.82'
$ws.Range("E34").Value = '  +2.47%  '
$ws.Range("D35").Value = '0.153'
$ws.Range("E35").Value = '  -5.80%  '
$ws.Range("E37").Value = '  +1.35%  '
$ws.Range("E38").Value = '  +2.39%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '177.54'
$ws.Range("E40").Value = '  +3.92%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  -1.82%  '
$ws.Range("E43").Value = '  +2.55%  '
$ws.Range("D45").Value = '1.94'
$ws.Range("E45").Value = '  +0.68%  '
$ws.Range("E46").Value = '  +2.19%  '
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").Value = '25.55'
$ws.Range("E48").Value = '  -2.32%  '
$ws.Range("E49").Value = '  +3.10%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '7.15'
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '23.51'
$ws.Range("E51").Value = '  +9.21%  '
